$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.394
$ws.Range("D5").Value = 0.481
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.531
$ws.Range("G5").Value = 0.497
$ws.Range("H5").Value = 0.518

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.394
$ws.Range("D7").Value = 0.481

# Row 8: BERT-base
$ws.Range("C8").Value = 0.383
$ws.Range("D8").Value = 0.585
$ws.Range("E8").Value = 0.6
$ws.Range("F8").Value = 0.645
$ws.Range("G8").Value = 0.618
$ws.Range("H8").Value = 0.637

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.367
$ws.Range("C9").Value = 0.55
$ws.Range("D9").Value = 0.65
$ws.Range("E9").Value = 0.664
$ws.Range("F9").Value = 0.701
$ws.Range("G9").Value = 0.656
$ws.Range("H9").Value = 0.666
